$wb = $excel.ActiveWorkbook

# Sheet "OFF": Week 15 logged -> update Road ("R") row (row 3)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 233
$wsOff.Range("C3").Value = 168
$wsOff.Range("D3").Value = 57
$wsOff.Range("G3").Value = 1

# Sheet "DEF": Week 16 simulated -> update Road ("R") row (row 3)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 220
$wsDef.Range("C3").Value = 148
$wsDef.Range("D3").Value = 59
$wsDef.Range("E3").Value = 26
